# This script applies the numeric corrections described in the commit
# "Update Name of Algo" to the KNN imputation result sheet.
# It updates 45 individual data cells (columns A, C, D) across rows 11-104
# of Sheet1 with their revised imputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("C11").Value = -12.288
$ws.Range("A12").Value = -21.488
$ws.Range("D14").Value = -8.17
$ws.Range("D19").Value = -8.034000000000001
$ws.Range("C23").Value = -12.394
$ws.Range("D24").Value = -7.344999999999999
$ws.Range("A27").Value = -21.8
$ws.Range("C28").Value = -12.992
$ws.Range("A32").Value = -22.02
$ws.Range("C32").Value = -13.305
$ws.Range("C34").Value = -12.101
$ws.Range("A36").Value = -20.464
$ws.Range("A38").Value = -20.03
$ws.Range("D38").Value = -8.184000000000001
$ws.Range("D41").Value = -8.15
$ws.Range("C42").Value = -12.334
$ws.Range("A46").Value = -21.733
$ws.Range("C49").Value = -13.035
$ws.Range("D52").Value = -7.679
$ws.Range("A54").Value = -21.809
$ws.Range("C54").Value = -13.054
$ws.Range("A55").Value = -21.961
$ws.Range("A56").Value = -21.912
$ws.Range("A67").Value = -21.577
$ws.Range("A69").Value = -21.47
$ws.Range("A72").Value = -21.624
$ws.Range("D72").Value = -7.464
$ws.Range("C78").Value = -12.21
$ws.Range("D78").Value = -7.444
$ws.Range("C80").Value = -12.695
$ws.Range("A83").Value = -22.005
$ws.Range("D83").Value = -7.998
$ws.Range("D85").Value = -8.664
$ws.Range("A86").Value = -22.115
$ws.Range("D86").Value = -8.617000000000001
$ws.Range("D90").Value = -6.994
$ws.Range("A91").Value = -20.755
$ws.Range("A93").Value = -21.54
$ws.Range("D96").Value = -7.762
$ws.Range("C97").Value = -11.768
$ws.Range("A99").Value = -21.703
$ws.Range("C99").Value = -12.177
$ws.Range("C101").Value = -12.188
$ws.Range("D103").Value = -8.344999999999999
$ws.Range("A104").Value = -21.175
